$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices (1-based) corresponding to C, D, E, F, G, K, M, N
$colIndices = @(3, 4, 5, 6, 7, 11, 13, 14)

$data = @(
    @(4.914449318578265, 7.100104438707438, 8.165091212297671, 38.9514190627079, 3.682017807928372, 19.03460693263725, 18.00567071725736, 19.65768971935569),
    @(4.752255673893469, 7.121964693118553, 8.081535132564687, 38.52403263842844, 3.686761693892341, 18.55268199265469, 17.73178154158508, 19.7148672339262),
    @(4.651601371514226, 7.136610636595108, 8.032426641897647, 38.27211468291809, 3.68981977575885, 18.2570447793237, 17.56744551242554, 19.75195348333597),
    @(4.610388031975169, 7.142884963468027, 8.012984279516331, 38.1721937351088, 3.691102674434596, 18.13683585091227, 17.50152577365778, 19.76756345983696),
    @(4.603534966160645, 7.143945241846358, 8.009790792312213, 38.1557700034562, 3.691317920518565, 18.11689721293324, 17.49064561056381, 19.7701854866089),
    @(4.651046247981508, 7.13669401766916, 8.032162106100587, 38.27075590539871, 3.689836928541129, 18.25542225727295, 17.56655214003683, 19.75216199324771),
    @(4.858788218845021, 7.107386846320933, 8.135836487348143, 38.80193772214606, 3.683623438622588, 18.86850100915665, 17.91048960199604, 19.67699340485789),
    @(5.254787083759568, 7.059702728651099, 8.355719762121163, 39.9222957166028, 3.67258431125447, 20.0645012490568, 18.61155991818817, 19.54531813316702),
    @(5.535263364064002, 7.030744660616292, 8.526175168160142, 40.78690182654016, 3.665161653216535, 20.92883539400312, 19.13771634625898, 19.45820366925825),
    @(5.659952239134157, 7.018914202737861, 8.605381955070314, 41.18778077258131, 3.661931973608348, 21.31690014385144, 19.37840122898501, 19.42067266508785),
    @(5.706705846388811, 7.01462944393506, 8.635591454214916, 41.34054563946083, 3.660729929125624, 21.46296039725792, 19.4696457235692, 19.4067632583013),
    @(5.69665790976322, 7.015543528961217, 8.629076092329582, 41.30760429480942, 3.660987881073293, 21.43154568543418, 19.44999143415435, 19.40974541379369),
    @(5.663808219345851, 7.018557769233107, 8.607863125546986, 41.20033030655968, 3.661832661445105, 21.32893531612044, 19.38590642055254, 19.41952225530056),
    @(5.64362520887569, 7.020429559366187, 8.594896935022948, 41.13474321913709, 3.662352839263013, 21.26596303334438, 19.34666321662752, 19.42555031653546),
    @(5.527052208457688, 7.031545006220353, 8.521030298805336, 40.76084482443781, 3.665375662987359, 20.90335783808854, 19.12200616119904, 19.46069871514459),
    @(5.454760883565259, 7.038709355267242, 8.476125223876917, 40.53331772425025, 3.667267582804296, 20.67948487301068, 18.98446355929693, 19.48279918599932),
    @(5.41291017530197, 7.042956354121069, 8.450455258883679, 40.40316993129355, 3.668369604027085, 20.55024512924079, 18.90548344596778, 19.49570817742915),
    @(5.398695251903488, 7.044415935957569, 8.441791753116735, 40.35923160819848, 3.668745111282096, 20.50641027863199, 18.87876735701346, 19.50011280835022),
    @(5.462484780609214, 7.037933616111941, 8.48088924033487, 40.55746479719154, 3.667064753605466, 20.70336674025631, 18.99909232706831, 19.48042611671392),
    @(5.673469887919447, 7.017667099795053, 8.614088228450962, 41.23181422021906, 3.661583961281129, 21.35909983303815, 19.40472766502171, 19.41664233371596),
    @(5.808642144333948, 7.005560410781427, 8.702386967215521, 41.67808111918557, 3.658124071704139, 21.78240024073411, 19.67039079067134, 19.37672127781429),
    @(5.736760875965902, 7.011917089186394, 8.655154263250592, 41.43943522971099, 3.659959559176659, 21.557005932813, 19.52857929221392, 19.39786598608197),
    @(5.458993704808391, 7.038283928879813, 8.4787349701783, 40.54654583835682, 3.667156408035009, 20.69257139388833, 18.99247835856621, 19.48149834869705),
    @(5.149262828184354, 7.071543639829549, 8.294580961918214, 39.61143775489654, 3.675449142938704, 19.7427497506019, 18.41958420199095, 19.57925205876161)
)

$startRow = 2
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $colIndices.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $colIndices[$c]).Value = $rowValues[$c]
    }
}

Write-Output "Updated loading_percent values for rows 2-25"